# Insert a new weekly price-report row for "Acelga" (Femacal de La Calera)
# above the existing row 549, shifting all subsequent rows (549-637) down by
# one, and fill the newly inserted row with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 549 - everything below (rows 549..637)
# shifts down to (550..638).
$ws.Rows.Item(549).Insert()

# Populate the newly inserted row 549 with the new data point.
$ws.Cells.Item(549, 1).Value  = 3
$ws.Cells.Item(549, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(549, 3).Value  = "Coquimbo"
$ws.Cells.Item(549, 4).Value  = 45218
$ws.Cells.Item(549, 5).Value  = 5
$ws.Cells.Item(549, 6).Value  = 100112009
$ws.Cells.Item(549, 7).Value  = "Acelga"
$ws.Cells.Item(549, 8).Value  = "Sin especificar"
$ws.Cells.Item(549, 9).Value  = "Primera"
$ws.Cells.Item(549, 10).Value = 220
$ws.Cells.Item(549, 11).Value = 3000
$ws.Cells.Item(549, 12).Value = 3500
$ws.Cells.Item(549, 13).Value = 3273
$ws.Cells.Item(549, 14).Value = "`$/docena de atados (6 kilos)"
$ws.Cells.Item(549, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(549, 16).Value = 546
$ws.Cells.Item(549, 17).Value = 6
$ws.Cells.Item(549, 18).Value = "Hortaliza"
